$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.420.35"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.826.47"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.86"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4586"
$ws.Range("E7").Value = "  -1.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3825"
$ws.Range("E8").Value = "  -1.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.15"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07900"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9713"
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.01"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.837.07"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.041"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.34"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06641"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.19"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.006"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.402.99"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.324"
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.309"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.042.21"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.01"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.37"
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.063"
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.235"
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.10"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9457"
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09303"
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.587"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.232"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.314"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05946"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02191"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.161"
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.006"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5761"
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1837"
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.01"
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.267"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5472"
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.98"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.867"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06634"
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "109.92"
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.005"
$ws.Range("E51").Value = "  -0.49%  "
